# Shift the Date values in A43:A436 down by one row (new[r] = old[r-1]),
# matching the OOXML diff: rows 43-436 each take on the date that
# previously occupied the row above them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "20220624", "20220625", "20220626", "20220627", "20220629", "20220701", "20220702", "20220703",
    "20220704", "20220705", "20220706", "20220707", "20220708", "20220709", "20220710", "20220711",
    "20220712", "20220713", "20220714", "20220715", "20220716", "20220717", "20220718", "20220719",
    "20220720", "20220721", "20220722", "20220723", "20220724", "20220726", "20220728", "20220729",
    "20220730", "20220731", "20220801", "20220802", "20220803", "20220804", "20220805", "20220806",
    "20220807", "20220808", "20220810", "20220811", "20220812", "20220813", "20220814", "20220815",
    "20220817", "20220818", "20220819", "20220820", "20220821", "20220822", "20220823", "20220824",
    "20220825", "20220826", "20220827", "20220828", "20220829", "20220830", "20220901", "20220902",
    "20220903", "20220904", "20220905", "20220906", "20220907", "20220908", "20220910", "20220911",
    "20220912", "20220913", "20220914", "20220915", "20220916", "20220917", "20220918", "20220919",
    "20220920", "20220921", "20220922", "20220923", "20220924", "20220925", "20220926", "20220927",
    "20220928", "20220929", "20221001", "20221002", "20221003", "20221004", "20221005", "20221006",
    "20221007", "20221008", "20221011", "20221012", "20221013", "20221014", "20221015", "20221016",
    "20221017", "20221018", "20221019", "20221020", "20221021", "20221022", "20221023", "20221024",
    "20221025", "20221026", "20221027", "20221028", "20221031", "20221030", "20221101", "20221102",
    "20221103", "20221104", "20221105", "20221107", "20221108", "20221109", "20221110", "20221111",
    "20221112", "20221113", "20221114", "20221115", "20221116", "20221117", "20221118", "20221119",
    "20221120", "20221121", "20221122", "20221123", "20221124", "20221125", "20221126", "20221127",
    "20221201", "20221202", "20221203", "20221204", "20221205", "20221206", "20221207", "20221208",
    "20221209", "20221210", "20221211", "20221212", "20221213", "20221214", "20221215", "20221216",
    "20221217", "20221218", "20221219", "20221220", "20221221", "20221222", "20221223", "20221224",
    "20221225", "20221226", "20221227", "20221228", "20221229", "20221230", "20230101", "20230102",
    "20230103", "20230104", "20230105", "20230106", "20230107", "20230108", "20230109", "20230110",
    "20230111", "20230112", "20230113", "20230114", "20230115", "20230116", "20230117", "20230118",
    "20230119", "20230120", "20230121", "20230122", "20230123", "20230124", "20230125", "20230126",
    "20230127", "20230128", "20230129", "20230130", "20230131", "20230201", "20230202", "20230203",
    "20230204", "20230206", "20230207", "20230208", "20230209", "20230210", "20230211", "20230212",
    "20230213", "20230214", "20230215", "20230217", "20230218", "20230219", "20230220", "20230221",
    "20230222", "20230223", "20230224", "20230225", "20230226", "20230227", "20230228", "20230301",
    "20230302", "20230303", "20230304", "20230305", "20230306", "20230307", "20230308", "20230309",
    "20230310", "20230311", "20230312", "20230313", "20230314", "20230315", "20230316", "20230317",
    "20230318", "20230319", "20230320", "20230321", "20230322", "20230323", "20230324", "20230325",
    "20230326", "20230327", "20230328", "20230329", "20230330", "20230331", "20230401", "20230402",
    "20230403", "20230404", "20230405", "20230406", "20230407", "20230408", "20230409", "20230410",
    "20230411", "20230412", "20230413", "20230414", "20230415", "20230416", "20230417", "20230418",
    "20230419", "20230420", "20230421", "20230422", "20230423", "20230424", "20230425", "20230426",
    "20230427", "20230428", "20230429", "20230430", "20230501", "20230502", "20230503", "20230504",
    "20230505", "20230506", "20230507", "20230508", "20230509", "20230510", "20230511", "20230512",
    "20230514", "20230515", "20230516", "20230517", "20230518", "20230519", "20230520", "20230521",
    "20230522", "20230523", "20230524", "20230525", "20230526", "20230527", "20230528", "20230529",
    "20230530", "20230531", "20230601", "20230602", "20230603", "20230604", "20230605", "20230606",
    "20230607", "20230608", "20230609", "20230610", "20230611", "20230612", "20230613", "20230614",
    "20230615", "20230616", "20230617", "20230618", "20230619", "20230620", "20230621", "20230622",
    "20230623", "20230624", "20230625", "20230626", "20230627", "20230628", "20230629", "20230630",
    "20230701", "20230702", "20230703", "20230704", "20230705", "20230706", "20230707", "20230708",
    "20230709", "20230710", "20230711", "20230712", "20230713", "20230714", "20230715", "20230716",
    "20230717", "20230718", "20230719", "20230720", "20230721", "20230722", "20230723", "20230725",
    "20230726", "20230727", "20230728", "20230729", "20190901", "20190902", "20190903", "20190906",
    "20190907", "20190908", "20190910", "20190911", "20190912", "20190913", "20190914", "20190915",
    "20190916", "20190917"
)

$startRow = 43
$endRow = 436
$count = $endRow - $startRow + 1

$targetRange = $ws.Range("A43:A436")
$targetRange.NumberFormat = "@"

$arr = New-Object "object[,]" $count,1
for ($i = 0; $i -lt $count; $i++) {
    $arr[$i, 0] = $newDates[$i]
}
$targetRange.Value = $arr
